$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.984.80"
$ws.Range("E2").Value = "  -1.99%  "

$ws.Range("D3").Value = "3.511.75"
$ws.Range("E3").Value = "  -0.83%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.88%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.635"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.47%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.634"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.24%  "

$ws.Range("E12").Value = "  -0.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.59%  "

$ws.Range("D14").Value = "4.077.76"
$ws.Range("E14").Value = "  -0.61%  "

$ws.Range("D15").Value = "3.508.90"
$ws.Range("E15").Value = "  -0.80%  "

$ws.Range("E16").Value = "  -0.42%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.87%  "

$ws.Range("D19").Value = "66.098.85"
$ws.Range("E19").Value = "  -1.82%  "

$ws.Range("E20").Value = "  -0.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "415.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.27%  "

$ws.Range("E22").Value = "  +4.69%  "

$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.69%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.81%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.96%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.18%  "

$ws.Range("E27").Value = "  -4.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "30.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "621.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.82%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.09%  "

$ws.Range("E32").Value = "  -2.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.111"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "59.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.87%  "

$ws.Range("E35").Value = "  +9.77%  "

$ws.Range("D36").Value = "0.0₃0807"
$ws.Range("E36").Value = "  -3.75%  "

$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.33%  "

$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.380"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.90%  "

$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "3.276.35"
$ws.Range("E40").Value = "  +7.99%  "

$ws.Range("E41").Value = "  -1.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.23%  "

$ws.Range("E43").Value = "  -4.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0419"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.79%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.75%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.67%  "

$ws.Range("E47").Value = "  -1.10%  "

$ws.Range("E48").Value = "  +0.64%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "138.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.15%  "
